$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update NoOfRoom value for the existing test case row (G2): "1 - One" -> "3 - Three"
$ws.Range("G2").Value = "3 - Three"

# The new value was entered with a distinct (non-theme) font color, as reflected
# by the new font/cell style introduced in the workbook (RGB 0x22,0x22,0x22).
$ws.Range("G2").Font.Color = 2236962

# Page orientation was set to portrait during this edit.
$ws.PageSetup.Orientation = 1

# Final selection left on the sheet after the edit.
$ws.Range("G5").Select()
